$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "DF" column header was removed from row 2/11/20 (the repeated table
# headers), shifting "FINAL Em." and "Regulation" one column to the left
# (F and G respectively) and leaving column H of those header rows blank.
foreach ($r in 2, 11, 20) {
    $ws.Range("F$r").Value = "FINAL Em."
    $ws.Range("G$r").Value = "Regulation"
    $ws.Range("H$r").Clear()
}
